# "terminamos de insertar datos"
# Finish filling in the ER-diagram data: update a couple of field
# definitions, add the missing FK/relationship cells, and tidy up the
# borders + column widths around the edited boxes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel border-edge constants (xlBordersIndex)
$xlEdgeLeft = 7
$xlEdgeTop = 8
$xlEdgeBottom = 9
$xlEdgeRight = 10
$xlContinuous = 1
$xlLineStyleNone = -4142
$xlPasteFormats = -4122

# ---------------------------------------------------------------
# Text / value updates
# ---------------------------------------------------------------

# Cursos.rut & Profesor.rut: INT [PK] -> VARCHAR(10) [PK]
$ws.Range("D5").Value = "rut VARCHAR(10) [PK]"
$ws.Range("H5").Value = "rut VARCHAR(10) [PK]"

# New relationship-multiplicity label between Profesor and Departamento
$ws.Range("I5").Value = "0..n           1"

# Prueba.profesor_calificador gains its type
$ws.Range("F7").Value = "profesor_calificador varchar(40)[FK]"

# Alumnos.id_curso gets a proper type and bracket position
$ws.Range("D8").Value = "id_curso INT NOT NULL [FK] "

# New Profesor.departamento [FK] field
$ws.Range("H8").Value = "departamento INT [FK]"

# New Alumnos.id_prueba [FK] field (brand new row in the Alumnos box,
# so first clone the plain-white / borderless formatting from the
# neighbouring untouched cell before typing into it)
$ws.Range("G9").Copy() | Out-Null
$ws.Range("D9").PasteSpecial($xlPasteFormats)
$ws.Range("D9").Value = "id_prueba INT NOT NULL [FK]"
$ws.Range("D9").Borders.Item($xlEdgeLeft).LineStyle = $xlContinuous
$ws.Range("D9").Borders.Item($xlEdgeRight).LineStyle = $xlContinuous
$ws.Range("D9").Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous

# ---------------------------------------------------------------
# Border touch-ups (box outlines shifting as rows were added)
# ---------------------------------------------------------------

# Alumnos / Profesor header cells: the box below them grew another
# row, so the header no longer closes off with a bottom rule.
$ws.Range("D4").Borders.Item($xlEdgeBottom).LineStyle = $xlLineStyleNone
$ws.Range("H4").Borders.Item($xlEdgeBottom).LineStyle = $xlLineStyleNone

# Profesor.rut now sits right under the header, so it needs its own
# top rule (matching the other "first field" cells).
$ws.Range("H5").Borders.Item($xlEdgeTop).LineStyle = $xlContinuous

# Multiplicity-label cells next to the boxes lose the inner edge that
# used to sit between them and the relationship label.
$ws.Range("C6").Borders.Item($xlEdgeRight).LineStyle = $xlLineStyleNone
$ws.Range("E6").Borders.Item($xlEdgeLeft).LineStyle = $xlLineStyleNone
$ws.Range("I6").Borders.Item($xlEdgeLeft).LineStyle = $xlLineStyleNone

# Apellido (Alumnos/Profesor) is no longer the last row of its box.
$ws.Range("H7").Borders.Item($xlEdgeBottom).LineStyle = $xlLineStyleNone
$ws.Range("D8").Borders.Item($xlEdgeBottom).LineStyle = $xlLineStyleNone

# New Profesor.departamento field closes off the bottom of the box.
$ws.Range("H8").Borders.Item($xlEdgeLeft).LineStyle = $xlContinuous
$ws.Range("H8").Borders.Item($xlEdgeRight).LineStyle = $xlContinuous
$ws.Range("H8").Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous

# ---------------------------------------------------------------
# Column widths (auto-fit no longer matches the new, longer text)
# ---------------------------------------------------------------
$ws.Columns.Item(6).ColumnWidth = 30.5
$ws.Columns.Item(9).ColumnWidth = 11

# ---------------------------------------------------------------
# Selection
# ---------------------------------------------------------------
$ws.Range("H8").Select()
